$wb = $excel.ActiveWorkbook

# Generate Report for handoff: mark the pending handoff rows as failed/ignored
# on both locale sheets (zh-cn, de-de). The Overview sheet is left untouched.
$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (B): handoff transform failed instead of ready
    $ws.Range("B2").Value = "Handoff transform failed"

    # Latest Handoff File column (C): drop the stale hyperlink + value entirely
    $ws.Range("C2").Hyperlinks.Delete()
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime column (D): reset to the "never happened" sentinel
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff Reason column (H): this row is now ignored, not included
    $ws.Range("H2").Value = "Ignored"
}
